# Applies the stat-table corrections described in the commit:
#   "Fixed README.md stats and docx preparation for all
#    Renaissance - JDK 17 - Z GC tests"
#
# The document is a single-column table where each row holds one
# (or, for three rows, several tab-separated) value(s). We address
# rows by their 1-based table-row index and overwrite the full cell
# range text, which keeps the existing run formatting (rFonts/sz).

$d = $word.ActiveDocument
$tbl = $d.Tables.Item(1)

# Row 1: 89.53 -> 0M
$tbl.Cell(1, 1).Range.Text = "0M"

# Row 2: 85.64 -> 0M
$tbl.Cell(2, 1).Range.Text = "0M"

# Row 3: 817 -> 0M
$tbl.Cell(3, 1).Range.Text = "0M"

# Row 4: 9061 -> 9064
$tbl.Cell(4, 1).Range.Text = "9064"

# Row 7: 0.04950 -> 0.04927
$tbl.Cell(7, 1).Range.Text = "0.04927"

# Row 8: 0.00578 -> 0.00569
$tbl.Cell(8, 1).Range.Text = "0.00569"

# Row 12: 85.53851 -> 85.64355
$tbl.Cell(12, 1).Range.Text = "85.64355"

# Row 44: collapse the multi-run/tab "1 ... 100.0" sequence down to "89.53"
$tbl.Cell(44, 1).Range.Text = "89.53"

# Row 45: collapse the multi-run/tab "1 ... 100.0" sequence down to "85.64"
$tbl.Cell(45, 1).Range.Text = "85.64"

# Row 46: collapse the multi-run/tab "1 ... 100.0" sequence down to "817"
$tbl.Cell(46, 1).Range.Text = "817"
